# Stroop task data file - add a new "neutral" condition block (rows 10-11)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: neutral stimulus "HHHHHH" shown in yellow ink -> correct key "n"
$ws.Cells.Item(10, 1).Value = "HHHHHH"
$ws.Cells.Item(10, 2).Value = "yellow"
$ws.Cells.Item(10, 4).Value = "n"

# Row 11: neutral stimulus "AAAAAA" shown in red ink -> correct key "m"
$ws.Cells.Item(11, 1).Value = "AAAAAA"
$ws.Cells.Item(11, 2).Value = "red"
$ws.Cells.Item(11, 4).Value = "m"

# Match the saved cursor position from the edited workbook
$ws.Range("C11").Select()
